{"js": "// Apply the LOM3056.docx content update:\n//  1. Ativa\u00e7\u00e3o date 2016 -> 2023\n//  2. Merge the two \"Objetivos\" sentences into one new paragraph of text\n//  3. Swap the responsible docente entry (Daniela -> Clodoaldo)\n//  4. Tweak the \"M\u00e9todo\" sentence (drop \"escritas\")\n//  5. Rewrite the \"Crit\u00e9rio\" grading formula/text\n//  6. Rewrite the \"Norma de recupera\u00e7\u00e3o\" text\n//  7. Update the \"Bibliografia\" references\n\nconst body = context.document.body;\n\n// 1. Ativa\u00e7\u00e3o date.\nconst ativacao = body.search(\"Ativa\u00e7\u00e3o: 01/01/2016\", { matchCase: true });\nativacao.load(\"items\");\nawait context.sync();\nif (ativacao.items.length > 0) {\n  ativacao.items[0].insertText(\"Ativa\u00e7\u00e3o: 01/01/2023\", Word.InsertLocation.replace);\n}\n\n// 2. Objetivos paragraph: the original paragraph holds a single run with two\n// <w:t> runs separated by a line break; replace the whole paragraph text\n// with the new single sentence (no line break).\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst oldObjetivosStart = \"- Apresentar aos alunos de Engenharia de Materiais\";\nlet objetivosParagraph = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(oldObjetivosStart) === 0) {\n    objetivosParagraph = paragraphs.items[i];\n    break;\n  }\n}\nif (objetivosParagraph) {\n  const newObjetivos =\n    \"Contextualizar os fundamentos de Qu\u00edmica Org\u00e2nica na \u00e1rea de Engenharia de Materiais, \" +\n    \"para permitir que os estudantes estejam aptos \u00e0 compreens\u00e3o das estruturas dos compostos \" +\n    \"org\u00e2nicos e sua influ\u00eancia nas propriedades dos materiais e dos principais mecanismos de \" +\n    \"rea\u00e7\u00e3o org\u00e2nica a serem utilizados na s\u00edntese e processamento de materiais pol\u00edmeros e \" +\n    \"outros materiais, como aqueles h\u00edbridos, por exemplo.\";\n  objetivosParagraph.insertText(newObjetivos, Word.InsertLocation.replace);\n}\n\n// 3. Responsible docente.\nconst docente = body.search(\"5840963 - Daniela Camargo Vernilli\", { matchCase: true });\ndocente.load(\"items\");\nawait context.sync();\nif (docente.items.length > 0) {\n  docente.items[0].insertText(\"5840897 - Clodoaldo Saron\", Word.InsertLocation.replace);\n}\n\n// 4. M\u00e9todo sentence.\nconst metodo = body.search(\"Avalia\u00e7\u00f5es escritas envolvendo o conte\u00fado da disciplina.\", { matchCase: true });\nmetodo.load(\"items\");\nawait context.sync();\nif (metodo.items.length > 0) {\n  metodo.items[0].insertText(\"Avalia\u00e7\u00f5es envolvendo o conte\u00fado da disciplina.\", Word.InsertLocation.replace);\n}\n\n// 5. Crit\u00e9rio text.\nconst criterioOld =\n  \"Duas avalia\u00e7\u00f5es no semestre (P1, P2). MS= (P1+P2)/2, onde: MS= m\u00e9dia do semestre. \" +\n  \"MS> ou = 5,0 = Aluno Aprovado MS< 3,0 = Aluno Reprovado 3,0 < ou = MS < 5,0 = Aluno de Recupera\u00e7\u00e3o.\";\nconst criterioNew =\n  \"Duas avalia\u00e7\u00f5es no semestre (P1, P2). MS= (2xP1+3xP2)/5, onde: MS= m\u00e9dia do semestre.\" +\n  \"MS> ou = 5,0: Aluno AprovadoMS< 3,0: Aluno Reprovado3,0 < ou = MS < 5,0: Aluno de Recupera\u00e7\u00e3o.\";\nconst criterio = body.search(criterioOld, { matchCase: true });\ncriterio.load(\"items\");\nawait context.sync();\nif (criterio.items.length > 0) {\n  criterio.items[0].insertText(criterioNew, Word.InsertLocation.replace);\n}\n\n// 6. Norma de recupera\u00e7\u00e3o text.\nconst normaOld =\n  \"Uma prova (PR), contendo todo o conte\u00fado da disciplina. O aluno ser\u00e1 aprovado se apresentar \" +\n  \"(m\u00e9dia final) MF > ou = 5,0. Onde: MF= MS+PR/2, onde: MS= m\u00e9dia do semestre e PR= prova de recupera\u00e7\u00e3o.\";\nconst normaNew =\n  \"Atividade avaliativa versando sobre o conte\u00fado da disciplina. O aluno ser\u00e1 aprovado se apresentar \" +\n  \"MF (m\u00e9dia final) > ou = 5,0. Onde: MF= (MS+PR)/2, onde: MS= m\u00e9dia do semestre e PR= prova de recupera\u00e7\u00e3o.\";\nconst norma = body.search(normaOld, { matchCase: true });\nnorma.load(\"items\");\nawait context.sync();\nif (norma.items.length > 0) {\n  norma.items[0].insertText(normaNew, Word.InsertLocation.replace);\n}\n\n// 7. Bibliografia text.\nconst bibliografiaOld =\n  \"McMURRY, J. Qu\u00edmica Org\u00e2nica. Rio de Janeiro: LTC Editora, 1997. - MORRISON, R.T. e BOYD, R.N. \" +\n  \"Qu\u00edmica Org\u00e2nica. 12\u00aa. Edi\u00e7\u00e3o. Lisboa: Fundac\u00e3o Calouste Gulbenkian, 1995. - SOLOMONS, T.W.G., \" +\n  \"FRYHLE, C.B. Qu\u00edmica Org\u00e2nica 1 e 2. 10\u00aa. Edi\u00e7\u00e3o, Rio de Janeiro: LTC Editora, 2012. - QUINO\u00c1, E. \" +\n  \"e RIGUERA, R. Quest\u00f5es e Exerc\u00edcios de Qu\u00edmica Org\u00e2nica. S\u00e3o Paulo: MAKRON Books, 1996.\";\nconst bibliografiaNew =\n  \"McMURRY, J. Qu\u00edmica Org\u00e2nica. 3\u00aa. Edi\u00e7\u00e3o. Editora Cengage Learning, 2016.- MORRISON, R.T. e BOYD, R.N. \" +\n  \"Qu\u00edmica Org\u00e2nica. 16\u00aa. Edi\u00e7\u00e3o. Lisboa: Fundac\u00e3o Calouste Gulbenkian, 2011.- SOLOMONS, T.W.G., FRYHLE, C.B. \" +\n  \"Qu\u00edmica Org\u00e2nica 1 e 2. 12\u00aa. Edi\u00e7\u00e3o, Rio de Janeiro: Gen/LTC Editora, 2018.\";\nconst bibliografia = body.search(bibliografiaOld, { matchCase: true });\nbibliografia.load(\"items\");\nawait context.sync();\nif (bibliografia.items.length > 0) {\n  bibliografia.items[0].insertText(bibliografiaNew, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Apply the LOM3056.docx content update:\n#  1. Ativa\u00e7\u00e3o date 2016 -> 2023\n#  2. Merge the two \"Objetivos\" sentences into one new paragraph of text\n#  3. Swap the responsible docente entry (Daniela -> Clodoaldo)\n#  4. Tweak the \"M\u00e9todo\" sentence (drop \"escritas\")\n#  5. Rewrite the \"Crit\u00e9rio\" grading formula/text\n#  6. Rewrite the \"Norma de recupera\u00e7\u00e3o\" text\n#  7. Update the \"Bibliografia\" references\n\n$d = $word.ActiveDocument\n\nfunction Replace-Text($old, $new) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.Text = $new\n    $find.Execute([ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]2, $true) | Out-Null\n}\n\n# 1. Ativa\u00e7\u00e3o date.\nReplace-Text \"Ativa\u00e7\u00e3o: 01/01/2016\" \"Ativa\u00e7\u00e3o: 01/01/2023\"\n\n# 2. Objetivos paragraph: originally a single run holding two <w:t> runs split\n# by a <w:br/>; replace the whole span (break included, via Chr(11)) with the\n# new single sentence containing no line break.\n$objetivosOld = \"- Apresentar aos alunos de Engenharia de Materiais os fundamentos de Qu\u00edmica Org\u00e2nica, para estarem aptos \u00e0 compreens\u00e3o das estruturas dos compostos org\u00e2nicos e dos principais mecanismos de rea\u00e7\u00e3o org\u00e2nica a serem utilizados principalmente na s\u00edntese e processamento de pol\u00edmeros.\" + [char]11 + \"- Ao final do curso, os alunos dever\u00e3o compreender e representar os mecanismos de processos org\u00e2nicos, compreender e representar  as equa\u00e7\u00f5es dos processos de obten\u00e7\u00e3o e propriedades qu\u00edmicas dos compostos org\u00e2nicos, entender a import\u00e2ncia das rea\u00e7\u00f5es org\u00e2nicas e dos processos de s\u00edntese org\u00e2nica e conhecer as principais t\u00e9cnicas de identifica\u00e7\u00e3o e caracteriza\u00e7\u00e3o de compostos org\u00e2nicos.\"\n$objetivosNew = \"Contextualizar os fundamentos de Qu\u00edmica Org\u00e2nica na \u00e1rea de Engenharia de Materiais, para permitir que os estudantes estejam aptos \u00e0 compreens\u00e3o das estruturas dos compostos org\u00e2nicos e sua influ\u00eancia nas propriedades dos materiais e dos principais mecanismos de rea\u00e7\u00e3o org\u00e2nica a serem utilizados na s\u00edntese e processamento de materiais pol\u00edmeros e outros materiais, como aqueles h\u00edbridos, por exemplo.\"\nReplace-Text $objetivosOld $objetivosNew\n\n# 3. Responsible docente.\nReplace-Text \"5840963 - Daniela Camargo Vernilli\" \"5840897 - Clodoaldo Saron\"\n\n# 4. M\u00e9todo sentence.\nReplace-Text \"Avalia\u00e7\u00f5es escritas envolvendo o conte\u00fado da disciplina.\" \"Avalia\u00e7\u00f5es envolvendo o conte\u00fado da disciplina.\"\n\n# 5. Crit\u00e9rio text.\n$criterioOld = \"Duas avalia\u00e7\u00f5es no semestre (P1, P2). MS= (P1+P2)/2, onde: MS= m\u00e9dia do semestre. MS> ou = 5,0 = Aluno Aprovado MS< 3,0 = Aluno Reprovado 3,0 < ou = MS < 5,0 = Aluno de Recupera\u00e7\u00e3o.\"\n$criterioNew = \"Duas avalia\u00e7\u00f5es no semestre (P1, P2). MS= (2xP1+3xP2)/5, onde: MS= m\u00e9dia do semestre.MS> ou = 5,0: Aluno AprovadoMS< 3,0: Aluno Reprovado3,0 < ou = MS < 5,0: Aluno de Recupera\u00e7\u00e3o.\"\nReplace-Text $criterioOld $criterioNew\n\n# 6. Norma de recupera\u00e7\u00e3o text.\n$normaOld = \"Uma prova (PR), contendo todo o conte\u00fado da disciplina. O aluno ser\u00e1 aprovado se apresentar (m\u00e9dia final) MF > ou = 5,0. Onde: MF= MS+PR/2, onde: MS= m\u00e9dia do semestre e PR= prova de recupera\u00e7\u00e3o.\"\n$normaNew = \"Atividade avaliativa versando sobre o conte\u00fado da disciplina. O aluno ser\u00e1 aprovado se apresentar MF (m\u00e9dia final) > ou = 5,0. Onde: MF= (MS+PR)/2, onde: MS= m\u00e9dia do semestre e PR= prova de recupera\u00e7\u00e3o.\"\nReplace-Text $normaOld $normaNew\n\n# 7. Bibliografia text.\n$bibliografiaOld = \"McMURRY, J. Qu\u00edmica Org\u00e2nica. Rio de Janeiro: LTC Editora, 1997. - MORRISON, R.T. e BOYD, R.N. Qu\u00edmica Org\u00e2nica. 12\u00aa. Edi\u00e7\u00e3o. Lisboa: Fundac\u00e3o Calouste Gulbenkian, 1995. - SOLOMONS, T.W.G., FRYHLE, C.B. Qu\u00edmica Org\u00e2nica 1 e 2. 10\u00aa. Edi\u00e7\u00e3o, Rio de Janeiro: LTC Editora, 2012. - QUINO\u00c1, E. e RIGUERA, R. Quest\u00f5es e Exerc\u00edcios de Qu\u00edmica Org\u00e2nica. S\u00e3o Paulo: MAKRON Books, 1996.\"\n$bibliografiaNew = \"McMURRY, J. Qu\u00edmica Org\u00e2nica. 3\u00aa. Edi\u00e7\u00e3o. Editora Cengage Learning, 2016.- MORRISON, R.T. e BOYD, R.N. Qu\u00edmica Org\u00e2nica. 16\u00aa. Edi\u00e7\u00e3o. Lisboa: Fundac\u00e3o Calouste Gulbenkian, 2011.- SOLOMONS, T.W.G., FRYHLE, C.B. Qu\u00edmica Org\u00e2nica 1 e 2. 12\u00aa. Edi\u00e7\u00e3o, Rio de Janeiro: Gen/LTC Editora, 2018.\"\nReplace-Text $bibliografiaOld $bibliografiaNew\n"}
